$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 22 price contains a unicode subscript-3 character (U+2083);
# build it from UTF-8 bytes since literal unicode escapes are not supported here.
$sub3Bytes = [byte[]]@(0xE2, 0x82, 0x83)
$sub3 = [System.Text.Encoding]::UTF8.GetString($sub3Bytes)

# Helper: some "Price" values look like plain numbers (e.g. "1.00", "0.0900").
# Assigning them straight to .Value lets Excel auto-convert to a Number and
# drop the significant trailing zeros, so those are entered with a leading
# apostrophe (the normal Excel 'force text' gesture) to keep the exact text.

$ws.Range("D2").Value = "52.349.85"
$ws.Range("E2").Value = "  +1.23%  "

$ws.Range("D3").Value = "2.893.34"
$ws.Range("E3").Value = "  +3.93%  "

$ws.Range("D4").Value = "'1.00"
$ws.Range("E4").Value = "  -0.11%  "

$ws.Range("D5").Value = "'353.03"
$ws.Range("E5").Value = "  +0.07%  "

$ws.Range("D6").Value = "'112.95"
$ws.Range("E6").Value = "  +3.60%  "

$ws.Range("D7").Value = "'0.561"
$ws.Range("E7").Value = "  +1.54%  "

$ws.Range("D8").Value = "'0.999"
$ws.Range("E8").Value = "  -0.32%  "

$ws.Range("E9").Value = "  +3.23%  "

$ws.Range("D10").Value = "'40.69"
$ws.Range("E10").Value = "  +2.05%  "

$ws.Range("E11").Value = "  -0.59%  "

$ws.Range("E12").Value = "  +2.02%  "

$ws.Range("D13").Value = "'20.34"
$ws.Range("E13").Value = "  +0.71%  "

$ws.Range("E14").Value = "  +2.51%  "

$ws.Range("D15").Value = "3.334.28"
$ws.Range("E15").Value = "  +3.32%  "

$ws.Range("B16").Value = "Polygon"
$ws.Range("C16").Value = "https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic"
$ws.Range("D16").Value = "'0.997"
$ws.Range("E16").Value = "  +7.58%  "

$ws.Range("B17").Value = "WrappedEther"
$ws.Range("C17").Value = "https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth"
$ws.Range("D17").Value = "2.873.70"
$ws.Range("E17").Value = "  +3.06%  "

$ws.Range("D18").Value = "52.269.13"
$ws.Range("E18").Value = "  +1.02%  "

$ws.Range("D19").Value = "'3.41"
$ws.Range("E19").Value = "  +8.88%  "

$ws.Range("D20").Value = "'7.70"
$ws.Range("E20").Value = "  -0.57%  "

$ws.Range("D21").Value = "'13.67"
$ws.Range("E21").Value = "  +3.90%  "

$ws.Range("D22").Value = "'0.0" + $sub3 + "0981"
$ws.Range("E22").Value = "  +1.56%  "

$ws.Range("D23").Value = "'70.84"
$ws.Range("E23").Value = "  +1.32%  "

$ws.Range("D24").Value = "'270.89"
$ws.Range("E24").Value = "  +1.48%  "

$ws.Range("D25").Value = "'2.79"
$ws.Range("E25").Value = "  +2.00%  "

$ws.Range("E26").Value = "  +1.86%  "

$ws.Range("E27").Value = "  -0.03%  "

$ws.Range("E28").Value = "  +1.31%  "

$ws.Range("E29").Value = "  +3.87%  "

$ws.Range("D30").Value = "'38.88"
$ws.Range("E30").Value = "  +4.34%  "

$ws.Range("D31").Value = "'6.34"
$ws.Range("E31").Value = "  +2.46%  "

$ws.Range("D32").Value = "'52.74"
$ws.Range("E32").Value = "  +1.92%  "

$ws.Range("D33").Value = "'0.0456"
$ws.Range("E33").Value = "  +0.57%  "

$ws.Range("D34").Value = "'0.0900"
$ws.Range("E34").Value = "  +8.26%  "

$ws.Range("D35").Value = "'5.65"
$ws.Range("E35").Value = "  +1.78%  "

$ws.Range("B36").Value = "FirstDigitalUSD"
$ws.Range("C36").Value = "https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd"
$ws.Range("D36").Value = "'1.00"
$ws.Range("E36").Value = "  -0.10%  "

$ws.Range("B37").Value = "Toncoin"
$ws.Range("C37").Value = "https://coinranking.com/coin/67YlI0K1b+toncoin-ton"
$ws.Range("D37").Value = "'1.93"
$ws.Range("E37").Value = "  -13.50%  "

$ws.Range("B38").Value = "LidoDAOToken"
$ws.Range("C38").Value = "https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo"
$ws.Range("D38").Value = "'3.34"
$ws.Range("E38").Value = "  +6.32%  "

$ws.Range("B39").Value = "Celestia"
$ws.Range("C39").Value = "https://coinranking.com/coin/YQcD0lBl7+celestia-tia"
$ws.Range("D39").Value = "'19.08"
$ws.Range("E39").Value = "  +3.08%  "

$ws.Range("E40").Value = "  +3.91%  "

$ws.Range("D41").Value = "'2.60"
$ws.Range("E41").Value = "  +2.52%  "

$ws.Range("E42").Value = "  +1.87%  "

$ws.Range("D43").Value = "'22.75"
$ws.Range("E43").Value = "  +2.74%  "

$ws.Range("D44").Value = "'122.84"
$ws.Range("E44").Value = "  +1.98%  "

$ws.Range("E45").Value = "  +2.00%  "

$ws.Range("D46").Value = "'3.58"
$ws.Range("E46").Value = "  +8.34%  "

$ws.Range("D47").Value = "2.186.21"
$ws.Range("E47").Value = "  +3.01%  "

$ws.Range("E48").Value = "  +7.21%  "

$ws.Range("D49").Value = "'0.246"
$ws.Range("E49").Value = "  +17.19%  "

$ws.Range("D50").Value = "'0.968"
$ws.Range("E50").Value = "  +6.28%  "

$ws.Range("D51").Value = "'0.0321"
$ws.Range("E51").Value = "  +12.71%  "

